# product_role_mapping.xlsx — update Display Product rule table.
#
# The "ProductDetails"/"UserDetails" rule-table header row (row 7) gains
# bound-variable prefixes ("product:" / "user:"), the productType/userRole
# condition cells (row 8) become real equality snippets, the action cell
# (row 8, column D) calls setDisplayProduct() on the new "product" binding,
# the ruleset-level import cell (C1) is replaced with the package import,
# and the stray "Variables"/descriptive cells (D7, B25, C25) are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ruleset import (C1): "Some business rules" -> the displayproduct package
$ws.Range("C1").Value = "poc.decisiontable.fuse.brms.displayproduct"

# RuleTable bound-variable header (row 7): add "product:" / "user:" prefixes
$ws.Range("B7").Value = "product:ProductDetails"
$ws.Range("C7").Value = "user:UserDetails"

# D7 no longer carries a description - clear it (keep the cell's style)
$ws.Range("D7").ClearContents()

# CONDITION/ACTION snippet row (row 8)
$ws.Range("B8").Value = 'productType=="$param"'
$ws.Range("C8").Value = 'userRole=="$param"'
$ws.Range("D8").Value = 'product.setDisplayProduct("$param");'

# Trailing "Variables" footer row is removed
$ws.Range("B25").ClearContents()
$ws.Range("C25").ClearContents()

# Selection moves from the (now nonexistent) D27 to C9
[void]$ws.Range("C9").Select()
